$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1=14, Q1=15 with the same style/format as the rest of row 1 (bold/border/centered)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# For rows 2-25: swap I<->K and M<->O values, and add new P and Q columns with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new, 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new, 2
}
